# "Generate Report for Handoff"
# The f5c94e38-c36f-4d86-9044-b2205eee830b.md file moves from "In Translation"
# to "Ready for handoff" status, and its per-language handoff timestamps are
# refreshed to reflect the newly generated handoff report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for f5c94e38-...-830b.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row for f5c94e38-...-830b.md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-08 20:14:00"

# --- de-de sheet: row for f5c94e38-...-830b.md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-08 20:14:07"
